$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "29.825.15"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +2.94%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.868.26"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +2.11%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.05%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "246.71"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.99%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.7013"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +2.27%  "

$ws.Cells.Item(7, 5).Value = "  -0.05%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.07783"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +2.11%  "

$ws.Cells.Item(9, 5).Value = "  +2.39%  "

$ws.Cells.Item(10, 5).Value = "  +2.38%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.07852"
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "5.196"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +3.01%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "1.872.63"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.98%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "92.95"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +2.85%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.6972"
$c.Style = "Normal"

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "6.671"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +3.51%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "29.816.65"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +2.85%  "

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "0.000008417"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.74%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "244.44"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.72%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "2.114.32"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.94%  "

$ws.Cells.Item(21, 5).Value = "  +1.47%  "

$ws.Cells.Item(22, 5).Value = "  -0.03%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "7.669"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +3.12%  "

$ws.Cells.Item(24, 5).Value = "  +0.03%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "0.1517"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +3.13%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "8.989"
$c.Style = "Normal"

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "160.27"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.63%  "

$ws.Cells.Item(28, 5).Value = "  +1.74%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "1.544"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.67%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "4.302"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +2.35%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "4.244"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +1.94%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "1.204"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.06%  "

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.05107"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.22%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "0.7918"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +4.43%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "1.937"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +6.63%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.170"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +1.78%  "

$ws.Cells.Item(37, 5).Value = "  +0.36%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "1.339.19"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +9.86%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.01890"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +3.01%  "

$ws.Cells.Item(40, 5).Value = "  +1.71%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.9746"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +6.95%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "6.051"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +11.90%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "107.17"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.43%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "1.0000"
$c.Style = "Normal"

$ws.Cells.Item(45, 5).Value = "  +4.00%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "9.838"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +4.17%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "2.013.69"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.62%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "65.66"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +3.54%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "1.800"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +4.21%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.5205"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.58%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "7.052"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +2.13%  "
